# Update odds values on row 2 of Sheet1 to reflect latest FlashScore data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value  = 3.15
$ws.Range("I2").Value  = 3.2
$ws.Range("K2").Value  = 2.1
$ws.Range("L2").Value  = 3.75
$ws.Range("Q2").Value  = 1.75
$ws.Range("R2").Value  = 2
$ws.Range("S2").Value  = 1.38
$ws.Range("T2").Value  = 2.82
$ws.Range("U2").Value  = 1.57
$ws.Range("V2").Value  = 2.25
$ws.Range("W2").Value  = 9.25
$ws.Range("X2").Value  = 12.5
$ws.Range("Z2").Value  = 24
$ws.Range("AA2").Value = 16
$ws.Range("AB2").Value = 21
$ws.Range("AD2").Value = 6.2
$ws.Range("AE2").Value = 11.5
$ws.Range("AI2").Value = 18
$ws.Range("AK2").Value = 45
$ws.Range("AL2").Value = 27
$ws.Range("AP2").Value = 16.5
$ws.Range("AT2").Value = 2.82
$ws.Range("AU2").Value = 6.4
$ws.Range("AW2").Value = 5.3
$ws.Range("AX2").Value = 18
$ws.Range("AY2").Value = 23
$ws.Range("AZ2").Value = 90
